$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "02/16/1990"
$ws.Range("Q3").Value = 45
$ws.Range("Y3").Value = "LA-2A compressor, 1176 compressor, Pultec"
$ws.Range("AA3").Value = "Melodyne, Fabfilter Bundle, iZotope Bundle"
$ws.Range("AB3").Value = "Moog, Roland, Nord Lead"
$ws.Range("G3").Value = "theweeknd@gmail.com"

$ws.Columns.Item(25).ColumnWidth = 36.916666666666664
$ws.Columns.Item(27).ColumnWidth = 36.583333333333336
$ws.Columns.Item(28).ColumnWidth = 21.25

$ws.Range("Q4").Select() | Out-Null
